$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing image URL (row 2) and fill in the remaining image URLs
# (rows 3-11) that were previously empty, mirroring the new Picasso-backed
# RecyclerView test data image links.
$ws.Range("C2").Value  = "http://i.imgur.com/zGvSDQJ.jpg"
$ws.Range("C3").Value  = "http://i.imgur.com/PUf63rI.jpg"
$ws.Range("C4").Value  = "http://i.imgur.com/fGrAGrb.png"
$ws.Range("C5").Value  = "http://i.imgur.com/VemIfxQ.jpg"
$ws.Range("C6").Value  = "http://i.imgur.com/EV5VY5m.jpg"
$ws.Range("C7").Value  = "http://i.imgur.com/WIHKX5U.jpg"
$ws.Range("C8").Value  = "http://i.imgur.com/pKI3a3m.jpg"
$ws.Range("C9").Value  = "http://i.imgur.com/pS1Q1Qa.jpg"
$ws.Range("C10").Value = "http://i.imgur.com/i2Fo7PN.jpg"
$ws.Range("C11").Value = "http://i.imgur.com/nenXqc0.jpg"

# Reflect the new active selection on the sheet (as saved in the workbook)
$ws.Range("C11").Select()
